# Generate Report for handback
# Updates the zh-cn and de-de localization-status sheets to reflect that the
# two content files have been handed back (in sync with en-US): fills in the
# "Latest Target File" / "Latest Handback File" hyperlink columns and updates
# the "Status" and "Latest Handback DateTime" values.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet ----------------------------------------------------------
# The Overview sheet mirrors the same "Status" shared string as the language
# sheets, so it picks up the new text as well.
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B2").Value = $statusText
$ws.Range("C2").Value = $statusText
$ws.Range("B3").Value = $statusText
$ws.Range("C3").Value = $statusText

# --- zh-cn sheet -----------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2 : 6193f53b-4736-4904-8a02-2c67061d7905.md
$ws.Range("B2").Value = $statusText
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/835bcf27731cc1137c74e48fc10dbb2a299e99ba/e2e/6193f53b-4736-4904-8a02-2c67061d7905.md", "", "", "6193f53b-4736-4904-8a02-2c67061d7905.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c9efcb56780d9a31abb16edf9b5ece0f368d5395/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/6193f53b-4736-4904-8a02-2c67061d7905.c0bedefabbbeb3c9b90a251f52d7936eca51b0fa.zh-cn.xlf", "", "", "6193f53b-4736-4904-8a02-2c67061d7905.c0bedefabbbeb3c9b90a251f52d7936eca51b0fa.zh-cn.xlf")
$ws.Range("G2").Value = "2016-01-25 10:55:39"

# Row 3 : ce536fa2-b3da-4f14-a74d-08faf91cb8e1.md
$ws.Range("B3").Value = $statusText
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/835bcf27731cc1137c74e48fc10dbb2a299e99ba/e2e/ce536fa2-b3da-4f14-a74d-08faf91cb8e1.md", "", "", "ce536fa2-b3da-4f14-a74d-08faf91cb8e1.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c9efcb56780d9a31abb16edf9b5ece0f368d5395/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ce536fa2-b3da-4f14-a74d-08faf91cb8e1.62ff99438f17fec30df2b756940a540b316895ae.zh-cn.xlf", "", "", "ce536fa2-b3da-4f14-a74d-08faf91cb8e1.62ff99438f17fec30df2b756940a540b316895ae.zh-cn.xlf")
$ws.Range("G3").Value = "2016-01-25 10:55:39"

# --- de-de sheet -------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 2 : 6193f53b-4736-4904-8a02-2c67061d7905.md
$ws.Range("B2").Value = $statusText
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/835bcf27731cc1137c74e48fc10dbb2a299e99ba/e2e/6193f53b-4736-4904-8a02-2c67061d7905.md", "", "", "6193f53b-4736-4904-8a02-2c67061d7905.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c3939f4da70ac3721a0ad29522b823ac3f6c9059/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/6193f53b-4736-4904-8a02-2c67061d7905.c0bedefabbbeb3c9b90a251f52d7936eca51b0fa.de-de.xlf", "", "", "6193f53b-4736-4904-8a02-2c67061d7905.c0bedefabbbeb3c9b90a251f52d7936eca51b0fa.de-de.xlf")
$ws.Range("G2").Value = "2016-01-25 10:55:55"

# Row 3 : ce536fa2-b3da-4f14-a74d-08faf91cb8e1.md
$ws.Range("B3").Value = $statusText
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/835bcf27731cc1137c74e48fc10dbb2a299e99ba/e2e/ce536fa2-b3da-4f14-a74d-08faf91cb8e1.md", "", "", "ce536fa2-b3da-4f14-a74d-08faf91cb8e1.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c3939f4da70ac3721a0ad29522b823ac3f6c9059/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ce536fa2-b3da-4f14-a74d-08faf91cb8e1.62ff99438f17fec30df2b756940a540b316895ae.de-de.xlf", "", "", "ce536fa2-b3da-4f14-a74d-08faf91cb8e1.62ff99438f17fec30df2b756940a540b316895ae.de-de.xlf")
$ws.Range("G3").Value = "2016-01-25 10:55:55"

Write-Output "Handback report generated"
